$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 6
$ws.Range("E6").Value = 7
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 13

# Set the active cell / selection to A6
$ws.Range("A6").Select()
